# research_notes.docx edit: merge spell-checked runs (drop proofErr markers)
# and append new marginalization notes at the end of the document.

$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    param($doc, [int]$index, [string]$finalText)
    $p = $doc.Paragraphs($index)
    $full = $p.Range
    $full.MoveEnd(1, -1) | Out-Null
    $null = $full.Find.Execute($full.Text, $true, $false, $false, $false, $false, $true, 1, $false, $finalText, 2)
}

function Merge-ParagraphRunsLeadingProofErr {
    param($doc, [int]$index, [string]$finalText)
    # Paragraph's first run is immediately preceded by a <w:proofErr> marker
    # (spellStart) that sits outside any run, so a plain Find/Replace across
    # the paragraph leaves it orphaned. Prefix a throwaway formatted character
    # so the marker falls inside the replaced span, then replace the whole lot.
    $p = $doc.Paragraphs($index)
    $r = $p.Range
    $r.Collapse(1) | Out-Null
    $r.InsertBefore("Z")
    $zRange = $doc.Range($p.Range.Start, $p.Range.Start + 1)
    $zRange.LanguageID = "en-US"
    $full = $p.Range
    $full.MoveEnd(1, -1) | Out-Null
    $searchText = "Z" + $finalText
    $null = $full.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $finalText, 2)
}

function Add-BlankParagraph {
    param($doc)
    # New paragraph that ends up as just <w:p><w:pPr><w:rPr><w:lang .../></w:rPr></w:pPr></w:p>:
    # type a placeholder char (inheriting en-US formatting), tag it, then delete
    # only the character (not the mark) so the language survives on the pilcrow.
    $last = $doc.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $doc.Paragraphs.Last
    $p.Style = "Normal"
    $r = $p.Range
    $r.InsertBefore("X")
    $charRange = $doc.Range($p.Range.Start, $p.Range.Start + 1)
    $charRange.LanguageID = "en-US"
    $charRange.Delete()
    $p.Range.LanguageID = "en-US"
}

function Add-NormalParagraph {
    param($doc, [string]$text)
    $last = $doc.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $doc.Paragraphs.Last
    $p.Style = "Normal"
    $r = $p.Range
    $r.InsertBefore("X")
    $full = $p.Range
    $full.MoveEnd(1, -1) | Out-Null
    $null = $full.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

function Add-ListParagraph {
    param($doc, [string]$text)
    # Inherits ListParagraph / ilvl=0 / numId=1 from the preceding list item.
    $last = $doc.Paragraphs.Last
    $last.Range.InsertParagraphAfter()
    $p = $doc.Paragraphs.Last
    $r = $p.Range
    $r.InsertBefore("X")
    $full = $p.Range
    $full.MoveEnd(1, -1) | Out-Null
    $null = $full.Find.Execute("X", $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

# --- 1. Clean up the proof-read runs (drop proofErr, merge text back into one run) ---

Merge-ParagraphRunsLeadingProofErr $d 1 "Carcinosarc as a factor"
Merge-ParagraphRuns $d 2 "Treating carcinosarc as 1 and all others as 0 as an additional factor"
Merge-ParagraphRuns $d 3 "Indeed carcinosarc is associated with worse prognosis with log.estimate of 0.338, though not significant due to high sd (low #’s)"
Merge-ParagraphRuns $d 10 "Do analysis where each TIL variable is considered independently(?) in a Coxph model, see which ones are best correlated "
Merge-ParagraphRuns $d 11 "Could do this perhaps by omitting the clust part of the coxph model on ipython notebook"

# --- 2. Append the new marginalization notes at the end of the document ---

Add-BlankParagraph $d
Add-BlankParagraph $d
Add-NormalParagraph $d "Molecular subtypes (p53abn, etc.) => vs angiogenesis markers "
Add-BlankParagraph $d
Add-NormalParagraph $d "When we get the HGSC data overlay and stuff => look at the S-TIL group (which should correspond to C1) => does that correspond to increased angiogenesis"
Add-BlankParagraph $d
Add-NormalParagraph $d "HER2 scores vs VEGFR"
Add-BlankParagraph $d
Add-NormalParagraph $d "Look at clustering HGSC CN signature exposures versus the endometrial ca ones using the Brenton samples + Vancouver samples"
Add-ListParagraph $d "Then do the same with TIL"
Add-ListParagraph $d "This is for Blake’s question"
Add-BlankParagraph $d
Add-BlankParagraph $d
Add-NormalParagraph $d "Talk about Wee1 inhibitors in CCNE1 mutated tumours cause we may have significant associations with the CCNE1 data from Dawn’s group"

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
